# Generate Report for Archive
#
# The localization-status report was regenerated. The entry for
# 98f95ff2-076f-4953-af01-54b66533f775 now sorts ahead of the entry for
# 18a579a9-720c-44a5-b53b-bf1e68f81538 (rows 5 and 6 swap on every sheet),
# and the 98f95ff2 entry's status moved from "Ready for handoff" back to
# "In Translation".

$wb = $excel.ActiveWorkbook

function Set-RowValues($ws, [int]$row, $values) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}

function Set-HyperlinkDisplay($ws, [string]$addr, [string]$display) {
    foreach ($h in $ws.Hyperlinks) {
        $a = $h.Range.Address()
        if ($a -eq $addr) {
            $h.TextToDisplay = $display
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-RowValues $wsOverview 5 @{
    "A" = "98f95ff2-076f-4953-af01-54b66533f775.md"
    "B" = "In Translation"
    "C" = "In Translation"
    "D" = "2016-03-23 00:37:46"
}
Set-RowValues $wsOverview 6 @{
    "A" = "18a579a9-720c-44a5-b53b-bf1e68f81538.md"
    "B" = "Ready for handoff"
    "C" = "Ready for handoff"
    "D" = "2016-03-23 00:38:10"
}

Set-HyperlinkDisplay $wsOverview '$A$5' "98f95ff2-076f-4953-af01-54b66533f775.md"
Set-HyperlinkDisplay $wsOverview '$A$6' "18a579a9-720c-44a5-b53b-bf1e68f81538.md"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-RowValues $wsZhCn 5 @{
    "A" = "98f95ff2-076f-4953-af01-54b66533f775.md"
    "C" = "In Translation"
    "D" = "98f95ff2-076f-4953-af01-54b66533f775.23b7bba38dc2e8c95e4a18a532490ed529bc76dd.zh-cn.xlf"
    "E" = "2016-03-23 00:37:41"
}
Set-RowValues $wsZhCn 6 @{
    "A" = "18a579a9-720c-44a5-b53b-bf1e68f81538.md"
    "C" = "Ready for handoff"
    "D" = "18a579a9-720c-44a5-b53b-bf1e68f81538.aed67083193761e77e6649b6061212086505938c.zh-cn.xlf"
    "E" = "2016-03-23 00:38:07"
}

Set-HyperlinkDisplay $wsZhCn '$A$5' "98f95ff2-076f-4953-af01-54b66533f775.md"
Set-HyperlinkDisplay $wsZhCn '$D$5' "98f95ff2-076f-4953-af01-54b66533f775.23b7bba38dc2e8c95e4a18a532490ed529bc76dd.zh-cn.xlf"
Set-HyperlinkDisplay $wsZhCn '$A$6' "18a579a9-720c-44a5-b53b-bf1e68f81538.md"
Set-HyperlinkDisplay $wsZhCn '$D$6' "18a579a9-720c-44a5-b53b-bf1e68f81538.aed67083193761e77e6649b6061212086505938c.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-RowValues $wsDeDe 5 @{
    "A" = "98f95ff2-076f-4953-af01-54b66533f775.md"
    "C" = "In Translation"
    "D" = "98f95ff2-076f-4953-af01-54b66533f775.23b7bba38dc2e8c95e4a18a532490ed529bc76dd.de-de.xlf"
    "E" = "2016-03-23 00:37:46"
}
Set-RowValues $wsDeDe 6 @{
    "A" = "18a579a9-720c-44a5-b53b-bf1e68f81538.md"
    "C" = "Ready for handoff"
    "D" = "18a579a9-720c-44a5-b53b-bf1e68f81538.aed67083193761e77e6649b6061212086505938c.de-de.xlf"
    "E" = "2016-03-23 00:38:10"
}

Set-HyperlinkDisplay $wsDeDe '$A$5' "98f95ff2-076f-4953-af01-54b66533f775.md"
Set-HyperlinkDisplay $wsDeDe '$D$5' "98f95ff2-076f-4953-af01-54b66533f775.23b7bba38dc2e8c95e4a18a532490ed529bc76dd.de-de.xlf"
Set-HyperlinkDisplay $wsDeDe '$A$6' "18a579a9-720c-44a5-b53b-bf1e68f81538.md"
Set-HyperlinkDisplay $wsDeDe '$D$6' "18a579a9-720c-44a5-b53b-bf1e68f81538.aed67083193761e77e6649b6061212086505938c.de-de.xlf"

Write-Output "Report archive regenerated: rows 5/6 swapped on Overview, zh-cn, de-de."
